$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: Insert a new blank row at position 23.
#
# This shifts old rows 23..33 down to 24..34 (preserving their styles/values),
# leaving a brand new empty row 23 ready to receive the "2001" data. Doing the
# insert at row 23 (rather than row 2) means the "Fonte" hyperlink text cells
# in O2:O5 are NOT shifted - they stay bound to their original row numbers,
# exactly like in the target workbook.
# ---------------------------------------------------------------------------
$ws.Rows("23").Insert()

# ---------------------------------------------------------------------------
# Step 2: Re-populate A2:M23 with the updated per-year data. Conceptually each
# year's dataset now lives one row further down than before (a new "2022" row
# was added at the top), so we simply write out the final values directly.
# Row 24 already contains the correct "2000" row courtesy of the insert above.
# ---------------------------------------------------------------------------
$data = @(
    @(2021, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2020, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2019, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2018, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2017, 0.4319, 0.41700000000000004, 0.27690000000000003, 0.39890000000000003, 0.24659999999999999, 0.3232, 0.3003, 0.309, 0.2976, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2016, 0.4721, 0.37889999999999996, 0.3425, 0.4639, 0.37729999999999997, 0.4003, 0.45139999999999997, 0.4091, 0.5016999999999999, 0.40449999999999997, 0.4071, 0.3897),
    @(2015, 0.35209999999999997, 0.3346, 0.2634, 0.3765, 0.3542, 0.36219999999999997, 0.4283, 0.47759999999999997, 0.4337, 0.439, 0.426, 0.3766),
    @(2014, 0.29610000000000003, 0.3595, 0.3004, 0.2732, 0.2926, 0.3071, 0.2932, 0.3522, 0.3069, 0.33409999999999995, 0.3506, 0.295),
    @(2013, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999, 0.26749999999999996, 0.24659999999999999, 0.2545, 0.3388, 0.2673),
    @(2012, 0.3405, 0.3332, 0.24659999999999999, 0.3536, 0.2693, 0.2935, 0.24659999999999999, 0.261, 0.2589, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2011, 0.3875, 0.3183, 0.3875, 0.3681, 0.2836, 0.404, 0.3583, 0.3698, 0.4547, 0.3471, 0.3087, 0.31120000000000003),
    @(2010, 0.3, 0.24659999999999999, 0.24659999999999999, 0.326, 0.24659999999999999, 0.2977, 0.30560000000000004, 0.362, 0.3377, 0.3169, 0.2939, 0.2803),
    @(2009, 0.462, 0.43099999999999994, 0.2918, 0.3907, 0.2921, 0.2916, 0.31229999999999997, 0.3519, 0.2663, 0.24659999999999999, 0.24659999999999999, 0.24659999999999999),
    @(2008, 0.3107, 0.3478, 0.27090000000000003, 0.2876, 0.3423, 0.32039999999999996, 0.36150000000000004, 0.4384, 0.4044, 0.44409999999999994, 0.49779999999999996, 0.40879999999999994),
    @(2007, 0.3992, 0.466, 0.3189, 0.4346, 0.3741, 0.41590000000000005, 0.3422, 0.3938, 0.3935, 0.2819, 0.36110000000000003, 0.30569999999999997),
    @(2006, 0.47400000000000003, 0.4797, 0.31930000000000003, 0.4544, 0.3323, 0.43579999999999997, 0.44079999999999997, 0.4221, 0.4908, 0.39909999999999995, 0.4345, 0.3751),
    @(2005, 0.48719999999999997, 0.43499999999999994, 0.34299999999999997, 0.5107, 0.44739999999999996, 0.4999, 0.5466000000000001, 0.5407000000000001, 0.594, 0.5109, 0.45710000000000006, 0.44),
    @(2004, 0.43689999999999996, 0.3749, 0.2925, 0.4248, 0.3342, 0.40159999999999996, 0.42310000000000003, 0.44229999999999997, 0.4476, 0.41980000000000006, 0.35760000000000003, 0.36150000000000004),
    @(2003, 0.6084, 0.7355999999999999, 0.6592, 0.6257, 0.666, 0.7127, 0.6642, 0.7944, 0.6514, 0.5838, 0.5687, 0.4246),
    @(2002, 0.44539999999999996, 0.5063, 0.364, 0.4228, 0.4829, 0.4573, 0.4052, 0.5127999999999999, 0.4953, 0.44260000000000005, 0.5241, 0.5115999999999999),
    @(2001, 0.3459, 0.3838, 0.2835, 0.4194, 0.40159999999999996, 0.42969999999999997, 0.3927, 0.49129999999999996, 0.5910000000000001, 0.4097, 0.5386, 0.4399)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $targetRow = $i + 3
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $values[$j]
    }
}

# ---------------------------------------------------------------------------
# Step 3: Row 2 becomes the new "2022" summary row. Only the January (column B)
# figure is populated; the rest of the monthly columns (C:M) are fully cleared
# out (not just emptied of content) so no leftover cell definitions remain.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 2022
$ws.Cells.Item(2, 2).Value = 0.29549999999999998
$ws.Range("C2:M2").Clear()
